# Drugs@FDA: remove submissions.submission_property_type.id row
# (internal primary key without a meaning) - shifts rows 44-47 up to 43-46.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43 holds: submission_property_type | id | string | The id of the submission property type.
$ws.Rows.Item(43).Delete()
